$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (volume number + week-covering dates) ---
$ws.Range("A8").Value = "Volume 30   Number  40"
$ws.Range("C9").Value = "Report Covering the Week  10/2/2023  Through  10/8/2023"

# --- Weekly crime-stat table updates (rows 14-30) ---
$ws.Range("C14").NumberFormat = '#,##0'
$ws.Range("C14").Value = 1
$ws.Range("F14").NumberFormat = '#,##0'
$ws.Range("F14").Value = 1
$ws.Range("I14").Value = 4
$ws.Range("K14").Value = -20
$ws.Range("L14").Value = -55.555555555555
$ws.Range("M14").Value = -75
$ws.Range("N14").Value = -81.818181818181
$ws.Range("D15").NumberFormat = '#,##0'
$ws.Range("D15").Value = 2
$ws.Range("E15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E15").Value = -100
$ws.Range("G15").Value = 2
$ws.Range("J15").Value = 18
$ws.Range("K15").Value = -5.555555555555
$ws.Range("N15").Value = -72.580645161290
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 0
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = -33.333333333333
$ws.Range("I16").Value = 125
$ws.Range("J16").Value = 130
$ws.Range("K16").Value = -3.846153846153
$ws.Range("L16").Value = -5.303030303030
$ws.Range("M16").Value = -52.651515151515
$ws.Range("N16").Value = -86.353711790393
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 60
$ws.Range("F17").Value = 19
$ws.Range("H17").Value = -40.625
$ws.Range("I17").Value = 238
$ws.Range("J17").Value = 265
$ws.Range("K17").Value = -10.188679245283
$ws.Range("L17").Value = -6.666666666666
$ws.Range("M17").Value = -6.299212598425
$ws.Range("N17").Value = -65.853658536585
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -75
$ws.Range("F18").Value = 5
$ws.Range("G18").Value = 21
$ws.Range("H18").Value = -76.190476190476
$ws.Range("I18").Value = 116
$ws.Range("J18").Value = 155
$ws.Range("K18").Value = -25.161290322580
$ws.Range("L18").Value = -10.077519379845
$ws.Range("M18").Value = -34.831460674157
$ws.Range("N18").Value = -81.229773462783
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = 11.111111111111
$ws.Range("F19").Value = 29
$ws.Range("G19").Value = 44
$ws.Range("H19").Value = -34.090909090909
$ws.Range("I19").Value = 258
$ws.Range("J19").Value = 333
$ws.Range("K19").Value = -22.522522522522
$ws.Range("L19").Value = -8.510638297872
$ws.Range("M19").Value = 1.976284584980
$ws.Range("N19").Value = -9.473684210526
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "0"
$ws.Range("D14").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 9
$ws.Range("G20").Value = 11
$ws.Range("H20").Value = -18.181818181818
$ws.Range("J20").Value = 107
$ws.Range("K20").Value = -22.429906542056
$ws.Range("L20").Value = 7.792207792207
$ws.Range("M20").Value = 18.571428571428
$ws.Range("N20").Value = -82.264957264957
$ws.Range("C21").Value = 21
$ws.Range("D21").Value = 22
$ws.Range("E21").Value = -4.545454545454
$ws.Range("F21").Value = 71
$ws.Range("G21").Value = 122
$ws.Range("H21").Value = -41.803278688524
$ws.Range("I21").Value = 841
$ws.Range("J21").Value = 1013
$ws.Range("K21").Value = -16.979269496544
$ws.Range("L21").Value = -5.823068309070
$ws.Range("M21").Value = -20.132953466286
$ws.Range("N21").Value = -72.588005215123
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0"
$ws.Range("D14").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("G22").Value = 1
$ws.Range("M22").Value = -52.380952380952
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "0"
$ws.Range("D14").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("D23").NumberFormat = '#,##0'
$ws.Range("D23").Value = 1
$ws.Range("E23").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E23").Value = -100
$ws.Range("F23").Value = 7
$ws.Range("G23").Value = 9
$ws.Range("H23").Value = -22.222222222222
$ws.Range("J23").Value = 67
$ws.Range("K23").Value = 1.492537313432
$ws.Range("L23").Value = -20
$ws.Range("C24").Value = 13
$ws.Range("D24").Value = 27
$ws.Range("E24").Value = -51.851851851851
$ws.Range("F24").Value = 74
$ws.Range("G24").Value = 109
$ws.Range("H24").Value = -32.110091743119
$ws.Range("I24").Value = 671
$ws.Range("J24").Value = 685
$ws.Range("K24").Value = -2.043795620437
$ws.Range("L24").Value = 45.238095238095
$ws.Range("M24").Value = 7.877813504823
$ws.Range("C25").Value = 10
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = 150
$ws.Range("F25").Value = 29
$ws.Range("G25").Value = 34
$ws.Range("H25").Value = -14.705882352941
$ws.Range("I25").Value = 384
$ws.Range("J25").Value = 321
$ws.Range("K25").Value = 19.626168224299
$ws.Range("L25").Value = 43.283582089552
$ws.Range("M25").Value = -38.950715421303
$ws.Range("D26").NumberFormat = '#,##0'
$ws.Range("D26").Value = 2
$ws.Range("E26").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E26").Value = -100
$ws.Range("G26").Value = 2
$ws.Range("J26").Value = 24
$ws.Range("K26").Value = 0
$ws.Range("C27").NumberFormat = '#,##0'
$ws.Range("C27").Value = 2
$ws.Range("D27").NumberFormat = '#,##0'
$ws.Range("D27").Value = 2
$ws.Range("E27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E27").Value = 0
$ws.Range("G27").NumberFormat = '#,##0'
$ws.Range("G27").Value = 2
$ws.Range("H27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H27").Value = 300
$ws.Range("I27").Value = 28
$ws.Range("J27").Value = 24
$ws.Range("K27").Value = 16.666666666666
$ws.Range("L27").Value = -30
$ws.Range("C28").NumberFormat = '#,##0'
$ws.Range("C28").Value = 1
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0"
$ws.Range("D14").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("F28").Value = 3
$ws.Range("H28").Value = -40
$ws.Range("I28").Value = 14
$ws.Range("K28").Value = -58.823529411764
$ws.Range("L28").Value = -58.823529411764
$ws.Range("M28").Value = -72
$ws.Range("N28").Value = -90.540540540540
$ws.Range("C29").NumberFormat = '#,##0'
$ws.Range("C29").Value = 1
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0"
$ws.Range("D14").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("F29").Value = 2
$ws.Range("H29").Value = -33.333333333333
$ws.Range("I29").Value = 11
$ws.Range("K29").Value = -57.692307692307
$ws.Range("L29").Value = -56
$ws.Range("M29").Value = -74.418604651162
$ws.Range("N29").Value = -91.666666666666
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0"
$ws.Range("D14").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("E30").PasteSpecial(-4122)
